$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "default"
$ws.Range("C3").Value = "Success"
$ws.Range("C4").Value = "Success"
